# Fruta / hortaliza, semanal
# Insert a new weekly record at row 48 (shifting existing rows 48:92 down to 49:93)
# and populate it with the latest observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(48).Insert()

$ws.Cells.Item(48, 1).Value = 9
$ws.Cells.Item(48, 2).Value = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(48, 3).Value = 'Metropolitana'
$ws.Cells.Item(48, 4).Value = 45240
$ws.Cells.Item(48, 5).Value = 13
$ws.Cells.Item(48, 6).Value = 100112010
$ws.Cells.Item(48, 7).Value = 'Achicoria'
$ws.Cells.Item(48, 8).Value = 'Sin especificar'
$ws.Cells.Item(48, 9).Value = 'Primera'
$ws.Cells.Item(48, 10).Value = 70
$ws.Cells.Item(48, 11).Value = 7000
$ws.Cells.Item(48, 12).Value = 8000
$ws.Cells.Item(48, 13).Value = 7500
$ws.Cells.Item(48, 14).Value = '$/caja 16 unidades'
$ws.Cells.Item(48, 15).Value = 'Provincia de Quillota'
$ws.Cells.Item(48, 16).Value = 469
$ws.Cells.Item(48, 17).Value = 16
$ws.Cells.Item(48, 18).Value = 'Hortaliza'
